$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 gets the data that was previously in row 7, and row 7 gets the
# data that was previously in row 6 (the two records traded places).

# --- Row 6 (new values, previously on row 7) ---
$ws.Range("A6").Value = 131046772
$ws.Range("B6").Value = 57884
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("M6").Value = "färska spår"
$ws.Range("Q6").Value = 401507
$ws.Range("R6").Value = 6818011
$ws.Range("Z6").Value = "15:17"
$ws.Range("AB6").Value = "15:17"
$ws.Range("AC6").Value = "Färska ringhack (tall)"

# --- Row 7 (new values, previously on row 6) ---
$ws.Range("A7").Value = 131046830
$ws.Range("B7").Value = 79243
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("M7").ClearContents()
$ws.Range("Q7").Value = 401538
$ws.Range("R7").Value = 6818009
$ws.Range("Z7").Value = "15:15"
$ws.Range("AB7").Value = "15:15"
$ws.Range("AC7").ClearContents()
